$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.521.22'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.619.82'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("E11").Value = '  +0.90%  '
$ws.Range("D12").Value = '1.846.09'
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").Value = '1.607.16'
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.23%  '
$ws.Range("E15").Value = '  -1.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '26.505.90'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +8.76%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  +2.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.94%  '
$ws.Range("D33").Value = '1.450.88'
$ws.Range("E33").Value = '  +8.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("E36").Value = '  +0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.560'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.90%  '
$ws.Range("E38").Value = '  +0.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.839'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.91%  '
$ws.Range("E40").Value = '  +2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.88%  '
$ws.Range("D43").Value = '1.758.07'
$ws.Range("E43").Value = '  +1.97%  '
$ws.Range("E44").Value = '  -0.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.77%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.917'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.33'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.22%  '
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("E49").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0965'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.91%  '
